$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "Options" column (currently column P / 16)
# so that "Antibiotics" (O) gets renamed, a new "Antibiotic_FinalConcentration"
# column is added right after it, and "Options" shifts one column to the right.
$ws.Columns.Item(16).Insert()

# Rename the existing "Antibiotics" header to "Antibiotic_name"
$ws.Cells.Item(1, 15).Value = "Antibiotic_name"

# New column header for the final concentration of the antibiotic
$ws.Cells.Item(1, 16).Value = "Antibiotic_FinalConcentration"

# Update the inducer A final-concentration values (column J, rows 3-8) so that
# the listing allows for floating point values instead of only integers
$newConc = "0_nM, 0.55nM, 0.10nM, 25nM, 50nM, 75_nM, 100nM"
$ws.Range("J3").Value = $newConc
$ws.Range("J4").Value = $newConc
$ws.Range("J5").Value = $newConc
$ws.Range("J6").Value = $newConc
$ws.Range("J7").Value = $newConc
$ws.Range("J8").Value = $newConc

# Match the new column widths used for the inserted columns
$ws.Columns.Item(15).ColumnWidth = 13.83
$ws.Columns.Item(16).ColumnWidth = 24.83

# Update view/selection state to match the saved workbook
$ws.Application.ActiveWindow.Zoom = 82
$ws.Range("J4:J8").Select()
